$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -0.9713005381550337
$ws.Cells.Item(2, 3).Value = 0.2910874087102975
$ws.Cells.Item(2, 4).Value = -0.1449385904233094
$ws.Cells.Item(2, 5).Value = 0.5154458957241742
$ws.Cells.Item(2, 6).Value = -0.1879827916858758
$ws.Cells.Item(2, 7).Value = 0.04615233206547309
$ws.Cells.Item(2, 8).Value = 0.08040507924932089
$ws.Cells.Item(2, 9).Value = 0.9057460109392793
$ws.Cells.Item(2, 10).Value = 0.249410584482785
$ws.Cells.Item(2, 11).Value = -0.4089899291462971
$ws.Cells.Item(3, 2).Value = 0.3128169072748576
$ws.Cells.Item(3, 3).Value = 0.7245973247384453
$ws.Cells.Item(3, 4).Value = -0.08695745817427486
$ws.Cells.Item(3, 5).Value = 0.1004139362101289
$ws.Cells.Item(3, 6).Value = 0.1144479125438225
$ws.Cells.Item(3, 7).Value = 0.9310525395423103
$ws.Cells.Item(3, 8).Value = 0.2709447666661666
$ws.Cells.Item(3, 9).Value = -0.3890834678507667
$ws.Cells.Item(3, 10).Value = 0.0092130989630414
$ws.Cells.Item(3, 11).Value = -0.05258007222182071
$ws.Cells.Item(4, 2).Value = -0.1030887313920102
$ws.Cells.Item(4, 3).Value = 0.07647302587122451
$ws.Cells.Item(4, 4).Value = 0.09012235911373839
$ws.Cells.Item(4, 5).Value = 0.9090979817469648
$ws.Cells.Item(4, 6).Value = 0.2510647147722038
$ws.Cells.Item(4, 7).Value = -0.4077104680353097
$ws.Cells.Item(4, 8).Value = -0.008762285969805494
$ws.Cells.Item(4, 9).Value = -0.07024161732427531
$ws.Cells.Item(4, 10).Value = -0.5044181462305073
$ws.Cells.Item(4, 11).Value = 0.4743807131573582
$ws.Cells.Item(5, 2).Value = 0.1720753300388297
$ws.Cells.Item(5, 3).Value = 0.9182953219789133
$ws.Cells.Item(5, 4).Value = 0.2313571522712326
$ws.Cells.Item(5, 5).Value = -0.4369487835323326
$ws.Cells.Item(5, 6).Value = -0.04090667078546828
$ws.Cells.Item(5, 7).Value = -0.1032375837743004
$ws.Cells.Item(5, 8).Value = -0.5376581911458389
$ws.Cells.Item(5, 9).Value = 0.4410716177777917
$ws.Cells.Item(5, 10).Value = 0.2139357511207785
$ws.Cells.Item(5, 11).Value = -0.4195295179412606
$ws.Cells.Item(6, 2).Value = 0.2025558615083408
$ws.Cells.Item(6, 3).Value = -0.4390725672109059
$ws.Cells.Item(6, 4).Value = -0.03362370600304576
$ws.Cells.Item(6, 5).Value = -0.09282107083105617
$ws.Cells.Item(6, 6).Value = -0.5262170457115171
$ws.Cells.Item(6, 7).Value = 0.4528464178518407
$ws.Cells.Item(6, 8).Value = 0.2258193719441708
$ws.Cells.Item(6, 9).Value = -0.4076102730996941
$ws.Cells.Item(6, 10).Value = 0.01979686535210479
$ws.Cells.Item(6, 11).Value = -0.1465366454910707
$ws.Cells.Item(7, 2).Value = -0.155505867098859
$ws.Cells.Item(7, 3).Value = -0.1846662397643801
$ws.Cells.Item(7, 4).Value = -0.6060897922829254
$ws.Cells.Item(7, 5).Value = 0.3779170903031727
$ws.Cells.Item(7, 6).Value = 0.1529367346250136
$ws.Cells.Item(7, 7).Value = -0.4796464385677615
$ws.Cells.Item(7, 8).Value = -0.0518896668375835
$ws.Cells.Item(7, 9).Value = -0.2180789225979328
$ws.Cells.Item(7, 10).Value = -0.09530727747379439
$ws.Cells.Item(7, 11).Value = -0.04804865355017168
$ws.Cells.Item(8, 2).Value = -0.5107750206255626
$ws.Cells.Item(8, 3).Value = 0.4561847949904575
$ws.Cells.Item(8, 4).Value = 0.2235272252108875
$ws.Cells.Item(8, 5).Value = -0.412585116272613
$ws.Cells.Item(8, 6).Value = 0.01354370435888463
$ws.Cells.Item(8, 7).Value = -0.1533995520910376
$ws.Cells.Item(8, 8).Value = -0.03097840982344718
$ws.Cells.Item(8, 9).Value = 0.01611673723937751
$ws.Cells.Item(8, 10).Value = -0.3805682542294263
$ws.Cells.Item(8, 11).Value = -0.1265478981343244
$ws.Cells.Item(9, 2).Value = 0.4836497670136274
$ws.Cells.Item(9, 3).Value = -0.2716791846901883
$ws.Cells.Item(9, 4).Value = 0.0989882439940924
$ws.Cells.Item(9, 5).Value = -0.09374781479368632
$ws.Cells.Item(9, 6).Value = 0.01667832048765883
$ws.Cells.Item(9, 7).Value = 0.058195260878178
$ws.Cells.Item(9, 8).Value = -0.3410837941741738
$ws.Cells.Item(9, 9).Value = -0.0882697463915933
$ws.Cells.Item(9, 10).Value = 0.1312785433800194
$ws.Cells.Item(9, 11).Value = -0.1047779946421779
$ws.Cells.Item(10, 2).Value = -0.09862000608775029
$ws.Cells.Item(10, 3).Value = -0.2177076384565043
$ws.Cells.Item(10, 4).Value = -0.07470130014448528
$ws.Cells.Item(10, 5).Value = -0.018762957560511
$ws.Cells.Item(10, 6).Value = -0.4116486652957548
$ws.Cells.Item(10, 7).Value = -0.1559957098169479
$ws.Cells.Item(10, 8).Value = 0.06481524009958373
$ws.Cells.Item(10, 9).Value = -0.1706787749263282
$ws.Cells.Item(10, 10).Value = -0.3495252635842865
$ws.Cells.Item(10, 11).Value = -0.3600316430428294
$ws.Cells.Item(11, 2).Value = -0.007299246851658558
$ws.Cells.Item(11, 3).Value = 0.05611473536190403
$ws.Cells.Item(11, 4).Value = -0.3334825685761003
$ws.Cells.Item(11, 5).Value = -0.07638392106091824
$ws.Cells.Item(11, 6).Value = 0.1450623997799764
$ws.Cells.Item(11, 7).Value = -0.09015246042256481
$ws.Cells.Item(11, 8).Value = -0.2688763367924936
$ws.Cells.Item(11, 9).Value = -0.2793288770644347
$ws.Cells.Item(11, 10).Value = -0.4870011301597575
$ws.Cells.Item(11, 11).Value = -0.2697718993153602
$ws.Cells.Item(12, 2).Value = -0.3458392694900738
$ws.Cells.Item(12, 3).Value = -0.09057491098672893
$ws.Cells.Item(12, 4).Value = 0.1301161455573347
$ws.Cells.Item(12, 5).Value = -0.1054089187162663
$ws.Cells.Item(12, 6).Value = -0.2842592456064009
$ws.Cells.Item(12, 7).Value = -0.2947628657683481
$ws.Cells.Item(12, 8).Value = -0.5024555220843632
$ws.Cells.Item(12, 9).Value = -0.285234325920563
$ws.Cells.Item(12, 10).Value = 0.02445755927687465
$ws.Cells.Item(12, 11).Value = -0.1298439175827339
$ws.Cells.Item(13, 2).Value = 0.3024188481467391
$ws.Cells.Item(13, 3).Value = -0.01091678595725426
$ws.Cells.Item(13, 4).Value = -0.2259338067192559
$ws.Cells.Item(13, 5).Value = -0.2532411366013347
$ws.Cells.Item(13, 6).Value = -0.4687411474418992
$ws.Cells.Item(13, 7).Value = -0.2551474109997061
$ws.Cells.Item(13, 8).Value = 0.05285908076345025
$ws.Cells.Item(13, 9).Value = -0.1022254650337468
$ws.Cells.Item(13, 10).Value = 0.5478383610817548
$ws.Cells.Item(13, 11).Value = 0.3315048434816775
$ws.Cells.Item(14, 2).Value = -0.3204423536017646
$ws.Cells.Item(14, 3).Value = -0.2897555356160152
$ws.Cells.Item(14, 4).Value = -0.4785434472588804
$ws.Cells.Item(14, 5).Value = -0.2526553077156036
$ws.Cells.Item(14, 6).Value = 0.06101006323956631
$ws.Cells.Item(14, 7).Value = -0.09146971300380768
$ws.Cells.Item(14, 8).Value = 0.5597931274032256
$ws.Cells.Item(14, 9).Value = 0.3440115546243878
$ws.Cells.Item(14, 10).Value = -0.06284176098887906
$ws.Cells.Item(14, 11).Value = 0.4362810619427306
$ws.Cells.Item(15, 2).Value = -0.3361387249989222
$ws.Cells.Item(15, 3).Value = -0.172079154149077
$ws.Cells.Item(15, 4).Value = 0.1140244002424372
$ws.Cells.Item(15, 5).Value = -0.05074182916399023
$ws.Cells.Item(15, 6).Value = 0.5950432373871204
$ws.Cells.Item(15, 7).Value = 0.376819120134967
$ws.Cells.Item(15, 8).Value = -0.03112348437038687
$ws.Cells.Item(15, 9).Value = 0.4675134828696286
$ws.Cells.Item(15, 10).Value = 0.308484786548516
$ws.Cells.Item(15, 11).Value = 0.3645688493654578
$ws.Cells.Item(16, 2).Value = 0.2919567931684033
$ws.Cells.Item(16, 3).Value = 0.03611576641084058
$ws.Cells.Item(16, 4).Value = 0.6379689631723193
$ws.Cells.Item(16, 5).Value = 0.3987277948059625
$ws.Cells.Item(16, 6).Value = -0.01926926647574417
$ws.Cells.Item(16, 7).Value = 0.4745601165931918
$ws.Cells.Item(16, 8).Value = 0.3132336451036749
$ws.Cells.Item(16, 9).Value = 0.3682199310910328
$ws.Cells.Item(16, 10).Value = 2.617463111368334
$ws.Cells.Item(16, 11).Value = 10.0863474470963
$ws.Cells.Item(17, 2).Value = 0.04668378654349858
$ws.Cells.Item(17, 3).Value = 0.646743751744091
$ws.Cells.Item(17, 4).Value = 0.4065301988496979
$ws.Cells.Item(17, 5).Value = -0.01195920564123576
$ws.Cells.Item(17, 6).Value = 0.4816209642468198
$ws.Cells.Item(17, 7).Value = 0.3201688216047165
$ws.Cells.Item(17, 8).Value = 0.3750919343242142
$ws.Cells.Item(17, 9).Value = 2.624303448075727
$ws.Cells.Item(17, 10).Value = 10.09317195097463
$ws.Cells.Item(17, 11).Value = -8.078224169346534
$ws.Cells.Item(18, 2).Value = 0.5397693479284105
$ws.Cells.Item(18, 3).Value = 0.3412705759114291
$ws.Cells.Item(18, 4).Value = -0.05754980417421507
$ws.Cells.Item(18, 5).Value = 0.4452902877305661
$ws.Cells.Item(18, 6).Value = 0.2881987268708067
$ws.Cells.Item(18, 7).Value = 0.3451757067707001
$ws.Cells.Item(18, 8).Value = 2.595354809876611
$ws.Cells.Item(18, 9).Value = 10.06467924472723
$ws.Cells.Item(18, 10).Value = -8.106501994311596
$ws.Cells.Item(18, 11).Value = 0.03320682975976369
$ws.Cells.Item(19, 2).Value = 0.3751437736559251
$ws.Cells.Item(19, 3).Value = -0.0547966439174179
$ws.Cells.Item(19, 4).Value = 0.4337013568798115
$ws.Cells.Item(19, 5).Value = 0.2700773128138926
$ws.Cells.Item(19, 6).Value = 0.3240755934079451
$ws.Cells.Item(19, 7).Value = 2.572896473190138
$ws.Cells.Item(19, 8).Value = 10.04160153156291
$ws.Cells.Item(19, 9).Value = -8.129862178460501
$ws.Cells.Item(19, 10).Value = 0.009717811984482905
$ws.Cells.Item(19, 11).Value = 2.157246767248676
$ws.Cells.Item(20, 2).Value = -0.312072073423668
$ws.Cells.Item(20, 3).Value = 0.3024220853976149
$ws.Cells.Item(20, 4).Value = 0.1967216315401449
$ws.Cells.Item(20, 5).Value = 0.2768343187407715
$ws.Cells.Item(20, 6).Value = 2.5374495888114
$ws.Cells.Item(20, 7).Value = 10.01148064496908
$ws.Cells.Item(20, 8).Value = -8.15757796247812
$ws.Cells.Item(20, 9).Value = -0.01691188077547812
$ws.Cells.Item(20, 10).Value = 2.13110753017287
$ws.Cells.Item(20, 11).Value = -1.308719802433826
$ws.Cells.Item(21, 2).Value = 0.259157971428892
$ws.Cells.Item(21, 3).Value = 0.16838023551327
$ws.Cells.Item(21, 4).Value = 0.2521582932286491
$ws.Cells.Item(21, 5).Value = 2.513162955834568
$ws.Cells.Item(21, 6).Value = 9.986857348713633
$ws.Cells.Item(21, 7).Value = -8.182561135418815
$ws.Cells.Item(21, 8).Value = -0.04214184030638374
$ws.Cells.Item(21, 9).Value = 2.105731958361166
$ws.Cells.Item(21, 10).Value = -1.334174982466894
$ws.Cells.Item(21, 11).Value = -1.395299246927819
$ws.Cells.Item(22, 2).Value = 0.2785713907394387
$ws.Cells.Item(22, 3).Value = 0.3082809271729959
$ws.Cells.Item(22, 4).Value = 2.547244771768479
$ws.Cells.Item(22, 5).Value = 10.01195673258608
$ws.Cells.Item(22, 6).Value = -8.161122484703073
$ws.Cells.Item(22, 7).Value = -0.02219513091605813
$ws.Cells.Item(22, 8).Value = 2.12507060979506
$ws.Cells.Item(22, 9).Value = -1.315084157547472
$ws.Cells.Item(22, 10).Value = -1.376309431020366
$ws.Cells.Item(22, 11).Value = 0.665125583809529
$ws.Cells.Item(23, 2).Value = 0.1514308419078968
$ws.Cells.Item(23, 3).Value = 2.452077041002839
$ws.Cells.Item(23, 4).Value = 9.957529601920154
$ws.Cells.Item(23, 5).Value = -8.197676472114839
$ws.Cells.Item(23, 6).Value = -0.04986031939925989
$ws.Cells.Item(23, 7).Value = 2.101663470123121
$ws.Cells.Item(23, 8).Value = -1.336421077423499
$ws.Cells.Item(23, 9).Value = -1.396642028765399
$ws.Cells.Item(23, 10).Value = 0.6452816047775823
$ws.Cells.Item(23, 11).Value = 0.06291644308519029
$ws.Cells.Item(24, 2).Value = 2.403213427390482
$ws.Cells.Item(24, 3).Value = 9.929760751643483
$ws.Cells.Item(24, 4).Value = -8.214896317152331
$ws.Cells.Item(24, 5).Value = -0.06195543535285469
$ws.Cells.Item(24, 6).Value = 2.092079342639022
$ws.Cells.Item(24, 7).Value = -1.344778008240451
$ws.Cells.Item(24, 8).Value = -1.404398746163505
$ws.Cells.Item(24, 9).Value = 0.6378183817016827
$ws.Cells.Item(24, 10).Value = 0.05559674216991822
$ws.Cells.Item(24, 11).Value = 0.1514544743675568
